# Commiting changes for CCDI test cases failed in jenkins
#
# The three SQL queries stored in column B (Participants/Samples/Files tab
# rows) are updated to cap their result sets with "LIMIT 100" so the CCDI
# test queries stop timing out in Jenkins.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 = ParticipantsTab query, B4 = SamplesTab query, B5 = FilesTab query
$participantQuery = $ws.Range("B2").Value2
$sampleQuery      = $ws.Range("B4").Value2
$fileQuery        = $ws.Range("B5").Value2

# Both the participant and sample queries already end in a newline before
# the trailing semicolon, so just swap the ";" for "LIMIT 100;".
$participantQuery = $participantQuery.Substring(0, $participantQuery.Length - 1) + "LIMIT 100;"
$sampleQuery      = $sampleQuery.Substring(0, $sampleQuery.Length - 1) + "LIMIT 100;"

# The file query's trailing semicolon directly follows the WHERE clause, so
# a newline needs to be introduced before the new LIMIT clause.
$fileQuery = $fileQuery.Substring(0, $fileQuery.Length - 1) + "`nLIMIT 100;"

$ws.Range("B2").Value2 = $participantQuery
$ws.Range("B4").Value2 = $sampleQuery
$ws.Range("B5").Value2 = $fileQuery

# Reflect the author's new active selection/view position on the sheet.
[void]$ws.Range("B4").Select()
